# fix: repair wrong data
# Corrects the "FrontEnd Deployed" (column H) figures on the Books_Create
# sheet, which had been accidentally populated with the wrong dataset, and
# restores the workbook's view state (selected cells / active sheet) to
# where the author left off after making the fix.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Data repair: Books_Create!H2:H61 ("FrontEnd Deployed" column) held
#    stale/incorrect values; replace them with the corrected figures.
# ---------------------------------------------------------------------
$wsBooksCreate = $wb.Worksheets.Item("Books_Create")

$correctedH = @(
    1720, 1230, 1140, 1230, 1740, 1170, 1200, 2050, 1130, 1130,
    1090, 1120, 1070, 1170, 1230, 1130, 1090, 1340, 1190, 1120,
    1150, 1130, 1120, 1080, 1280, 1110, 1090, 1130, 1120, 1130,
    1080, 1130, 1230, 1080, 1130, 1120, 1130, 1130, 3580, 1230,
    1250, 1140, 1230, 1250, 1140, 1230, 1160, 1190, 1230, 1230,
    1230, 1160, 1160, 1200, 1230, 1630, 1240, 1240, 1230, 1150
)

$row = 2
foreach ($value in $correctedH) {
    $wsBooksCreate.Cells.Item($row, 8).Value = $value
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Restore each sheet's own cursor position / selection.
#    Selecting a range on a worksheet also activates that worksheet, so
#    the final selection made below (on Books_Create) leaves it as the
#    active / displayed tab, matching the saved workbook state.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Transactions_All").Range("A1").Select()
$wb.Worksheets.Item("Transactions_Id").Range("C30").Select()
$wb.Worksheets.Item("Customers_with_Profile").Range("H2").Select()
$wb.Worksheets.Item("Books_Update").Range("A2").Select()
$wb.Worksheets.Item("Books_Delete").Range("A2").Select()
$wsBooksCreate.Range("J6").Select()
